# EDTEch.pptx edit: insert a new "Feasibility Analysis" section-header slide
# right after slide 4 ("Elixir"), and tidy up the run structure of the
# "Advantages" slide's bullet list (merge previously-split runs back into a
# single run per bullet line).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Insert the new slide at position 5, using the same "1_Section header"
#    custom layout used by the neighbouring section-header slides
#    (Elixir / Working / Business Impact).
# ---------------------------------------------------------------------
$sectionHeaderLayout = $p.SlideMaster.CustomLayouts.Item(14)
$newSlide = $p.Slides.AddSlide(5, $sectionHeaderLayout)

# Placeholder 1: centered title ("ctrTitle")
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Feasibility Analysis"
$titleShape.Left = 47.68299212598425
$titleShape.Top = 34.50007874015748
$titleShape.Width = 492.31692913385825
$titleShape.Height = 66

# Placeholder 2: subtitle body text
$subShape = $newSlide.Shapes.Item(2)
$para1 = "Due to covid19 pandemic situation the governments across the world   looking for ways to shift education to online platforms due to the pandemic situation. As a result of this, many platforms have emerged which provide interaction between teachers and students. Though online teaching culture is gaining widespread attention, they are not infallible and the government is giving out surveys to validate the effectiveness of these tools. EDtech  possesses the capability of coupling with any of these platforms and providing an accurate analysis of student/teacher engagement."
$para3 = "education administration: It becomes difficult for the administration like schools, colleges,etc to have an unbiased feedback of the students for the faculty. This leads to incompetence, posing threat to the quality of education. Our product is capable of taking the attentiveness statistics of the whole class and using mathematical calculations to analyse the effectiveness of a faculty. For example, let" + [char]0x2019 + "s say we calculate the average attentiveness for the whole class and if it is low, we can conclude that either the faculty is not putting enough effort or his/her pedagogy is a bit screwed. Accordingly, further actions can be taken."
$para5 = "Teaching faculty: Virtually, it becomes impossible for the faculty members to track every student. Although it doesn" + [char]0x2019 + "t affect other students" + [char]0x2019 + " performance, it certainly results in casual behaviour and short attention span of students. Using this tool faculty members can see reports of individual students and their performance during the lecture. Consistently poor performing students can be identified and necessary actions maybe taken."

$lines = @($para1, "", $para3, "", $para5, "", "")
$subShape.TextFrame.TextRange.Text = [string]::Join([char]13, $lines)
$subShape.Left = 47.68299212598425
$subShape.Top = 100.5
$subShape.Width = 654.3169291338583
$subShape.Height = 270

# Placeholder 3: secondary title (idx 2) is left blank, matching the source.

# ---------------------------------------------------------------------
# 2. On the "Advantages" slide (now slide 16, pushed down by the insert
#    above), collapse each bullet's split runs back into a single run.
# ---------------------------------------------------------------------
$advSlide = $p.Slides.Item(16)
$bodyShape = $advSlide.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

$bodyRange.Characters(1, 15).Text = "Fully Automatic"
$bodyRange.Characters(18, 34).Text = "Enhance Creativity & Visualization"
$bodyRange.Characters(54, 37).Text = "Improve Student-Teacher Collaboration"
$bodyRange.Characters(93, 18).Text = "Incessant Learning"
